# Apply updated "F" column (attendance/visit-count style metric) values
# across the four worksheets, per the target diff.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 58
$ws.Range("F5").Value  = 538
$ws.Range("F7").Value  = 1376
$ws.Range("F9").Value  = 971
$ws.Range("F14").Value = 3964
$ws.Range("F17").Value = 3062
$ws.Range("F18").Value = 831
$ws.Range("F19").Value = 146
$ws.Range("F21").Value = 141
$ws.Range("F22").Value = 2065
$ws.Range("F24").Value = 1922
$ws.Range("F28").Value = 8536
$ws.Range("F29").Value = 5748
$ws.Range("F38").Value = 41
$ws.Range("F41").Value = 4647
$ws.Range("F43").Value = 864

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F18").Value = 1130

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 365

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value  = 365
$ws.Range("F7").Value  = 58
$ws.Range("F8").Value  = 538
$ws.Range("F9").Value  = 1376
$ws.Range("F11").Value = 971
$ws.Range("F15").Value = 3964
$ws.Range("F17").Value = 3062
$ws.Range("F18").Value = 831
$ws.Range("F19").Value = 146
$ws.Range("F21").Value = 2065
$ws.Range("F27").Value = 1922
$ws.Range("F31").Value = 8536
$ws.Range("F32").Value = 5748
$ws.Range("F43").Value = 4647
$ws.Range("F44").Value = 864
